# ---------------------------------------------------------------------------
# Updates cached numeric values on the "Leve Profits" sheets (ARM, BSM, CRP,
# CUL, GSM, LTW, WVR) to reflect refreshed marketboard prices pulled in by the
# scheduled Sheets runner. Every touched cell is a plain literal value (the
# workbook has no formulas), so this script writes each updated price/profit
# figure directly via the Excel object model.
#
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ,
#          J=currentAveragePriceHQ, K=LevePriceNQ, L=LevePriceHQ,
#          M=LeveProfitNQ, N=LeveProfitHQ
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ARM")
# Row 32: H32=11088.82, I32=9613.775, K32=9613.775, M32=-9326.775
$ws.Cells.Item(32, 8).Value2 = 11088.82
$ws.Cells.Item(32, 9).Value2 = 9613.775
$ws.Cells.Item(32, 11).Value2 = 9613.775
$ws.Cells.Item(32, 13).Value2 = -9326.775
# Row 61: H61=3021.6667, I61=1658.55, K61=1658.55, M61=-1446.55
$ws.Cells.Item(61, 8).Value2 = 3021.6667
$ws.Cells.Item(61, 9).Value2 = 1658.55
$ws.Cells.Item(61, 11).Value2 = 1658.55
$ws.Cells.Item(61, 13).Value2 = -1446.55
# Row 74: H74=1199.375, I74=1299.8182, J74=978.4, K74=1299.8182, L74=978.4, M74=-425.8181999999999, N74=-2726.4
$ws.Cells.Item(74, 8).Value2 = 1199.375
$ws.Cells.Item(74, 9).Value2 = 1299.8182
$ws.Cells.Item(74, 10).Value2 = 978.4
$ws.Cells.Item(74, 11).Value2 = 1299.8182
$ws.Cells.Item(74, 12).Value2 = 978.4
$ws.Cells.Item(74, 13).Value2 = -425.8181999999999
$ws.Cells.Item(74, 14).Value2 = -2726.4
# Row 76: H76=79624.75, J76=78166.336, L76=78166.336, N76=-78842.336
$ws.Cells.Item(76, 8).Value2 = 79624.75
$ws.Cells.Item(76, 10).Value2 = 78166.336
$ws.Cells.Item(76, 12).Value2 = 78166.336
$ws.Cells.Item(76, 14).Value2 = -78842.336
# Row 77: H77=1199.375, I77=1299.8182, J77=978.4, K77=6499.090999999999, L77=4892, M77=-2131.090999999999, N77=-13628
$ws.Cells.Item(77, 8).Value2 = 1199.375
$ws.Cells.Item(77, 9).Value2 = 1299.8182
$ws.Cells.Item(77, 10).Value2 = 978.4
$ws.Cells.Item(77, 11).Value2 = 6499.090999999999
$ws.Cells.Item(77, 12).Value2 = 4892
$ws.Cells.Item(77, 13).Value2 = -2131.090999999999
$ws.Cells.Item(77, 14).Value2 = -13628
# Row 79: H79=79624.75, J79=78166.336, L79=78166.336, N79=-80506.336
$ws.Cells.Item(79, 8).Value2 = 79624.75
$ws.Cells.Item(79, 10).Value2 = 78166.336
$ws.Cells.Item(79, 12).Value2 = 78166.336
$ws.Cells.Item(79, 14).Value2 = -80506.336
# Row 132: H132=950, I132=950, K132=2850, M132=-320
$ws.Cells.Item(132, 8).Value2 = 950
$ws.Cells.Item(132, 9).Value2 = 950
$ws.Cells.Item(132, 11).Value2 = 2850
$ws.Cells.Item(132, 13).Value2 = -320
# Row 136: H136=3021.6667, I136=1658.55, K136=4975.65, M136=-2425.65
$ws.Cells.Item(136, 8).Value2 = 3021.6667
$ws.Cells.Item(136, 9).Value2 = 1658.55
$ws.Cells.Item(136, 11).Value2 = 4975.65
$ws.Cells.Item(136, 13).Value2 = -2425.65

$ws = $wb.Worksheets.Item("BSM")
# Row 86: H86=3277.1614, I86=1909.5, J86=5170.846, K86=1909.5, L86=5170.846, M86=-786.5, N86=-7416.846
$ws.Cells.Item(86, 8).Value2 = 3277.1614
$ws.Cells.Item(86, 9).Value2 = 1909.5
$ws.Cells.Item(86, 10).Value2 = 5170.846
$ws.Cells.Item(86, 11).Value2 = 1909.5
$ws.Cells.Item(86, 12).Value2 = 5170.846
$ws.Cells.Item(86, 13).Value2 = -786.5
$ws.Cells.Item(86, 14).Value2 = -7416.846
# Row 89: H89=3277.1614, I89=1909.5, J89=5170.846, K89=9547.5, L89=25854.23, M89=-3931.5, N89=-37086.23
$ws.Cells.Item(89, 8).Value2 = 3277.1614
$ws.Cells.Item(89, 9).Value2 = 1909.5
$ws.Cells.Item(89, 10).Value2 = 5170.846
$ws.Cells.Item(89, 11).Value2 = 9547.5
$ws.Cells.Item(89, 12).Value2 = 25854.23
$ws.Cells.Item(89, 13).Value2 = -3931.5
$ws.Cells.Item(89, 14).Value2 = -37086.23
# Row 94: H94=4935.0527, I94=4063.9167, J94=6428.4287, K94=4063.9167, L94=6428.4287, M94=-3612.9167, N94=-7330.4287
$ws.Cells.Item(94, 8).Value2 = 4935.0527
$ws.Cells.Item(94, 9).Value2 = 4063.9167
$ws.Cells.Item(94, 10).Value2 = 6428.4287
$ws.Cells.Item(94, 11).Value2 = 4063.9167
$ws.Cells.Item(94, 12).Value2 = 6428.4287
$ws.Cells.Item(94, 13).Value2 = -3612.9167
$ws.Cells.Item(94, 14).Value2 = -7330.4287

$ws = $wb.Worksheets.Item("CRP")
# Row 43: H43=12000, J43=12000, L43=12000, N43=-12368
$ws.Cells.Item(43, 8).Value2 = 12000
$ws.Cells.Item(43, 10).Value2 = 12000
$ws.Cells.Item(43, 12).Value2 = 12000
$ws.Cells.Item(43, 14).Value2 = -12368
# Row 58: H58=1550.7931, I58=767.35297, K58=767.35297, M58=-564.35297
$ws.Cells.Item(58, 8).Value2 = 1550.7931
$ws.Cells.Item(58, 9).Value2 = 767.35297
$ws.Cells.Item(58, 11).Value2 = 767.35297
$ws.Cells.Item(58, 13).Value2 = -564.35297
# Row 69: H69=18000, I69=18000, K69=18000, M69=-17251
$ws.Cells.Item(69, 8).Value2 = 18000
$ws.Cells.Item(69, 9).Value2 = 18000
$ws.Cells.Item(69, 11).Value2 = 18000
$ws.Cells.Item(69, 13).Value2 = -17251
# Row 72: H72=18000, I72=18000, K72=54000, M72=-50256
$ws.Cells.Item(72, 8).Value2 = 18000
$ws.Cells.Item(72, 9).Value2 = 18000
$ws.Cells.Item(72, 11).Value2 = 54000
$ws.Cells.Item(72, 13).Value2 = -50256
# Row 74: H74=61314.4, J74=61314.4, L74=61314.4, N74=-63062.4
$ws.Cells.Item(74, 8).Value2 = 61314.4
$ws.Cells.Item(74, 10).Value2 = 61314.4
$ws.Cells.Item(74, 12).Value2 = 61314.4
$ws.Cells.Item(74, 14).Value2 = -63062.4
# Row 77: H77=61314.4, J77=61314.4, L77=183943.2, N77=-192679.2
$ws.Cells.Item(77, 8).Value2 = 61314.4
$ws.Cells.Item(77, 10).Value2 = 61314.4
$ws.Cells.Item(77, 12).Value2 = 183943.2
$ws.Cells.Item(77, 14).Value2 = -192679.2
# Row 95: H95=12032.2, J95=12032.2, L95=12032.2, N95=-17524.2
$ws.Cells.Item(95, 8).Value2 = 12032.2
$ws.Cells.Item(95, 10).Value2 = 12032.2
$ws.Cells.Item(95, 12).Value2 = 12032.2
$ws.Cells.Item(95, 14).Value2 = -17524.2
# Row 101: H101=12000, J101=12000, L101=12000, N101=-18490
$ws.Cells.Item(101, 8).Value2 = 12000
$ws.Cells.Item(101, 10).Value2 = 12000
$ws.Cells.Item(101, 12).Value2 = 12000
$ws.Cells.Item(101, 14).Value2 = -18490
# Row 102: H102=30241, J102=30241, L102=30241, N102=-35109
$ws.Cells.Item(102, 8).Value2 = 30241
$ws.Cells.Item(102, 10).Value2 = 30241
$ws.Cells.Item(102, 12).Value2 = 30241
$ws.Cells.Item(102, 14).Value2 = -35109
# Row 132: H132=3980, I132=3980, J132=0, K132=11940, L132=0, N132=-9410
$ws.Cells.Item(132, 8).Value2 = 3980
$ws.Cells.Item(132, 9).Value2 = 3980
$ws.Cells.Item(132, 10).Value2 = 0
$ws.Cells.Item(132, 11).Value2 = 11940
$ws.Cells.Item(132, 12).Value2 = 0
$ws.Cells.Item(132, 14).Value2 = -9410
# Row 136: H136=1550.7931, I136=767.35297, K136=2302.05891, M136=247.9410899999998
$ws.Cells.Item(136, 8).Value2 = 1550.7931
$ws.Cells.Item(136, 9).Value2 = 767.35297
$ws.Cells.Item(136, 11).Value2 = 2302.05891
$ws.Cells.Item(136, 13).Value2 = 247.9410899999998
# Row 141: H141=112173.086, J141=112173.086, L141=112173.086, N141=-122533.086
$ws.Cells.Item(141, 8).Value2 = 112173.086
$ws.Cells.Item(141, 10).Value2 = 112173.086
$ws.Cells.Item(141, 12).Value2 = 112173.086
$ws.Cells.Item(141, 14).Value2 = -122533.086

$ws = $wb.Worksheets.Item("CUL")
# Row 9: H9=758.3333, I9=310, J9=3000, K9=930, L9=9000, M9=-706, N9=-9448
$ws.Cells.Item(9, 8).Value2 = 758.3333
$ws.Cells.Item(9, 9).Value2 = 310
$ws.Cells.Item(9, 10).Value2 = 3000
$ws.Cells.Item(9, 11).Value2 = 930
$ws.Cells.Item(9, 12).Value2 = 9000
$ws.Cells.Item(9, 13).Value2 = -706
$ws.Cells.Item(9, 14).Value2 = -9448
# Row 137: H137=1166.4286, J137=903.5714, L137=2710.7142, N137=-12910.7142
$ws.Cells.Item(137, 8).Value2 = 1166.4286
$ws.Cells.Item(137, 10).Value2 = 903.5714
$ws.Cells.Item(137, 12).Value2 = 2710.7142
$ws.Cells.Item(137, 14).Value2 = -12910.7142
# Row 139: H139=3042.7856, I139=2781.7273, J139=4000, K139=8345.1819, L139=12000, M139=-3205.1819, N139=-22280
$ws.Cells.Item(139, 8).Value2 = 3042.7856
$ws.Cells.Item(139, 9).Value2 = 2781.7273
$ws.Cells.Item(139, 10).Value2 = 4000
$ws.Cells.Item(139, 11).Value2 = 8345.1819
$ws.Cells.Item(139, 12).Value2 = 12000
$ws.Cells.Item(139, 13).Value2 = -3205.1819
$ws.Cells.Item(139, 14).Value2 = -22280

$ws = $wb.Worksheets.Item("GSM")
# Row 62: H62=40974.5, J62=40974.5, L62=40974.5, N62=-42346.5
$ws.Cells.Item(62, 8).Value2 = 40974.5
$ws.Cells.Item(62, 10).Value2 = 40974.5
$ws.Cells.Item(62, 12).Value2 = 40974.5
$ws.Cells.Item(62, 14).Value2 = -42346.5
# Row 65: H65=40974.5, J65=40974.5, L65=122923.5, N65=-129787.5
$ws.Cells.Item(65, 8).Value2 = 40974.5
$ws.Cells.Item(65, 10).Value2 = 40974.5
$ws.Cells.Item(65, 12).Value2 = 122923.5
$ws.Cells.Item(65, 14).Value2 = -129787.5
# Row 97: H97=495.2, I97=499.84616, J97=465, K97=499.84616, L97=465, M97=-3.846159999999998, N97=-1457
$ws.Cells.Item(97, 8).Value2 = 495.2
$ws.Cells.Item(97, 9).Value2 = 499.84616
$ws.Cells.Item(97, 10).Value2 = 465
$ws.Cells.Item(97, 11).Value2 = 499.84616
$ws.Cells.Item(97, 12).Value2 = 465
$ws.Cells.Item(97, 13).Value2 = -3.846159999999998
$ws.Cells.Item(97, 14).Value2 = -1457
# Row 102: H102=1048.8235, I102=867.6923, K102=867.6923, M102=754.3077
$ws.Cells.Item(102, 8).Value2 = 1048.8235
$ws.Cells.Item(102, 9).Value2 = 867.6923
$ws.Cells.Item(102, 11).Value2 = 867.6923
$ws.Cells.Item(102, 13).Value2 = 754.3077
# Row 122: H122=3344, I122=1441.6428, K122=4324.928400000001, M122=-1874.928400000001
$ws.Cells.Item(122, 8).Value2 = 3344
$ws.Cells.Item(122, 9).Value2 = 1441.6428
$ws.Cells.Item(122, 11).Value2 = 4324.928400000001
$ws.Cells.Item(122, 13).Value2 = -1874.928400000001
# Row 126: H126=4793.4443, I126=3410.25, K126=10230.75, M126=-7760.75
$ws.Cells.Item(126, 8).Value2 = 4793.4443
$ws.Cells.Item(126, 9).Value2 = 3410.25
$ws.Cells.Item(126, 11).Value2 = 10230.75
$ws.Cells.Item(126, 13).Value2 = -7760.75
# Row 132: H132=3128.5, I132=2166.6667, K132=6500.000100000001, M132=-3970.000100000001
$ws.Cells.Item(132, 8).Value2 = 3128.5
$ws.Cells.Item(132, 9).Value2 = 2166.6667
$ws.Cells.Item(132, 11).Value2 = 6500.000100000001
$ws.Cells.Item(132, 13).Value2 = -3970.000100000001
# Row 133: H133=70000, J133=70000, L133=70000, N133=-80120
$ws.Cells.Item(133, 8).Value2 = 70000
$ws.Cells.Item(133, 10).Value2 = 70000
$ws.Cells.Item(133, 12).Value2 = 70000
$ws.Cells.Item(133, 14).Value2 = -80120
# Row 136: H136=12320.667, J136=12320.667, L136=36962.001, N136=-42062.001
$ws.Cells.Item(136, 8).Value2 = 12320.667
$ws.Cells.Item(136, 10).Value2 = 12320.667
$ws.Cells.Item(136, 12).Value2 = 36962.001
$ws.Cells.Item(136, 14).Value2 = -42062.001

$ws = $wb.Worksheets.Item("LTW")
# Row 16: H16=455.55554, I16=429.5, K16=429.5, M16=-259.5
$ws.Cells.Item(16, 8).Value2 = 455.55554
$ws.Cells.Item(16, 9).Value2 = 429.5
$ws.Cells.Item(16, 11).Value2 = 429.5
$ws.Cells.Item(16, 13).Value2 = -259.5
# Row 40: H40=9030.174000000001, I40=9637.615, K40=9637.615, M40=-9501.615
$ws.Cells.Item(40, 8).Value2 = 9030.174000000001
$ws.Cells.Item(40, 9).Value2 = 9637.615
$ws.Cells.Item(40, 11).Value2 = 9637.615
$ws.Cells.Item(40, 13).Value2 = -9501.615
# Row 46: H46=1392.7142, I46=1674.75, J46=1016.6667, K46=1674.75, L46=1016.6667, M46=-1486.75, N46=-1392.6667
$ws.Cells.Item(46, 8).Value2 = 1392.7142
$ws.Cells.Item(46, 9).Value2 = 1674.75
$ws.Cells.Item(46, 10).Value2 = 1016.6667
$ws.Cells.Item(46, 11).Value2 = 1674.75
$ws.Cells.Item(46, 12).Value2 = 1016.6667
$ws.Cells.Item(46, 13).Value2 = -1486.75
$ws.Cells.Item(46, 14).Value2 = -1392.6667
# Row 93: H93=2374.2942, I93=2197.4167, K93=2197.4167, M93=-949.4167000000002
$ws.Cells.Item(93, 8).Value2 = 2374.2942
$ws.Cells.Item(93, 9).Value2 = 2197.4167
$ws.Cells.Item(93, 11).Value2 = 2197.4167
$ws.Cells.Item(93, 13).Value2 = -949.4167000000002
# Row 122: H122=4872.0835, J122=6522, L122=19566, N122=-24466
$ws.Cells.Item(122, 8).Value2 = 4872.0835
$ws.Cells.Item(122, 10).Value2 = 6522
$ws.Cells.Item(122, 12).Value2 = 19566
$ws.Cells.Item(122, 14).Value2 = -24466
# Row 132: H132=6139.8184, I132=2794.3333, J132=7394.375, K132=8382.999899999999, L132=22183.125, M132=-5852.999899999999, N132=-27243.125
$ws.Cells.Item(132, 8).Value2 = 6139.8184
$ws.Cells.Item(132, 9).Value2 = 2794.3333
$ws.Cells.Item(132, 10).Value2 = 7394.375
$ws.Cells.Item(132, 11).Value2 = 8382.999899999999
$ws.Cells.Item(132, 12).Value2 = 22183.125
$ws.Cells.Item(132, 13).Value2 = -5852.999899999999
$ws.Cells.Item(132, 14).Value2 = -27243.125
# Row 136: H136=4218.069, I136=2574.9524, J136=8531.25, K136=7724.8572, L136=25593.75, M136=-5174.8572, N136=-30693.75
$ws.Cells.Item(136, 8).Value2 = 4218.069
$ws.Cells.Item(136, 9).Value2 = 2574.9524
$ws.Cells.Item(136, 10).Value2 = 8531.25
$ws.Cells.Item(136, 11).Value2 = 7724.8572
$ws.Cells.Item(136, 12).Value2 = 25593.75
$ws.Cells.Item(136, 13).Value2 = -5174.8572
$ws.Cells.Item(136, 14).Value2 = -30693.75

$ws = $wb.Worksheets.Item("WVR")
# Row 4: H4=8568.857, I4=5000, J4=9163.666999999999, K4=5000, L4=9163.666999999999, M4=-4887, N4=-9389.666999999999
$ws.Cells.Item(4, 8).Value2 = 8568.857
$ws.Cells.Item(4, 9).Value2 = 5000
$ws.Cells.Item(4, 10).Value2 = 9163.666999999999
$ws.Cells.Item(4, 11).Value2 = 5000
$ws.Cells.Item(4, 12).Value2 = 9163.666999999999
$ws.Cells.Item(4, 13).Value2 = -4887
$ws.Cells.Item(4, 14).Value2 = -9389.666999999999
# Row 100: H100=481.7857, I100=508.34784, J100=359.6, K100=1016.69568, L100=719.2, M100=-475.69568, N100=-1801.2
$ws.Cells.Item(100, 8).Value2 = 481.7857
$ws.Cells.Item(100, 9).Value2 = 508.34784
$ws.Cells.Item(100, 10).Value2 = 359.6
$ws.Cells.Item(100, 11).Value2 = 1016.69568
$ws.Cells.Item(100, 12).Value2 = 719.2
$ws.Cells.Item(100, 13).Value2 = -475.69568
$ws.Cells.Item(100, 14).Value2 = -1801.2
# Row 132: H132=2943.923, I132=2327.1, K132=6981.299999999999, M132=-4451.299999999999
$ws.Cells.Item(132, 8).Value2 = 2943.923
$ws.Cells.Item(132, 9).Value2 = 2327.1
$ws.Cells.Item(132, 11).Value2 = 6981.299999999999
$ws.Cells.Item(132, 13).Value2 = -4451.299999999999
# Row 136: H136=1477.0615, I136=671.7959, K136=2015.3877
$ws.Cells.Item(136, 8).Value2 = 1477.0615
$ws.Cells.Item(136, 9).Value2 = 671.7959
$ws.Cells.Item(136, 11).Value2 = 2015.3877
